$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48: fix amplitude bug for NiBecCameraOdt trial ---
# K48: StrongLight -> TwoLevelWeakLight
$ws.Cells.Item(48, 11).Value = "TwoLevelWeakLight"
# L48 (AdCLim): 8 -> 4
$ws.Cells.Item(48, 12).Value = 4

# --- Row 64 (new): partialEvapDpartialevaptime ---
$ws.Cells.Item(64, 1).Value = "partialEvapDpartialevaptime"
$ws.Cells.Item(64, 2).Value = "An experiment at partial evaporation stage D."
$ws.Cells.Item(64, 3).Value = "TOP"
$ws.Cells.Item(64, 4).Value = "EvapDOdt1"
$ws.Cells.Item(64, 5).Value = "None"
$ws.Cells.Item(64, 6).Value = 4
$ws.Cells.Item(64, 7).Value = "partialevaptime"
$ws.Cells.Item(64, 8).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(64, 9).Value = "LSR"
$ws.Cells.Item(64, 10).Value = "LF"
$ws.Cells.Item(64, 11).Value = "RandomPolarization"
$ws.Cells.Item(64, 12).Value = 8
$ws.Cells.Item(64, 13).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(64, 14).Value = 30
$ws.Cells.Item(64, 15).Value = "LinearFit1D"

# --- Row 65 (new): HfBecCurveField ---
$ws.Cells.Item(65, 1).Value = "HfBecCurveField"
$ws.Cells.Item(65, 2).Value = "An experiment at the high-field BEC stage. Scan curveField."
$ws.Cells.Item(65, 3).Value = "TOP"
$ws.Cells.Item(65, 4).Value = "Bec"
$ws.Cells.Item(65, 5).Value = "None"
$ws.Cells.Item(65, 6).Value = 4
$ws.Cells.Item(65, 7).Value = "curveField"
$ws.Cells.Item(65, 8).Value = "DensityFit;AtomNumber"
$ws.Cells.Item(65, 9).Value = "LSR"
$ws.Cells.Item(65, 10).Value = "HF"
$ws.Cells.Item(65, 11).Value = "StrongLight"
$ws.Cells.Item(65, 12).Value = 8
$ws.Cells.Item(65, 13).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(65, 14).Value = 2
$ws.Cells.Item(65, 15).Value = "LinearFit1D"

# --- Column A got a touch wider after the edits (was 25.5703125) ---
$ws.Columns.Item(1).ColumnWidth = 25.65
